# This script applies the target edit:
#  - Renames sheet "事業投資" to "具有相當價值之財產" and replaces its data with a single
#    "珠寶" (jewelry) record.
#  - Adds three new sheets after it, in order: "保險" (insurance), "債權" (receivables/bonds),
#    and a brand-new "事業投資" (enterprise investment) sheet which receives the data that used
#    to live in the original "事業投資" sheet.
#  - Rewrites sheet "其他有價證券" (5th sheet) replacing its old/garbled rows with a single,
#    correctly structured record.

$wb = $excel.ActiveWorkbook

function Write-HeaderAndRow($ws, $headers, $values, $styleSrc) {
    $lastCol = $headers.Count + 1  # data starts at column B (col 2), A is reserved for index
    # Header row (row 1), starting at column B
    for ($i = 0; $i -lt $headers.Count; $i++) {
        $col = $i + 2
        $ws.Cells.Item(1, $col).Value = $headers[$i]
    }
    # Data row (row 2), starting at column B; column A holds the numeric index value (values[0])
    $ws.Cells.Item(2, 1).Value = $values[0]
    for ($i = 0; $i -lt $headers.Count; $i++) {
        $col = $i + 2
        $ws.Cells.Item(2, $col).Value = $values[$i + 1]
    }

    # Copy cell formatting from a known-good template sheet so styles match
    # (column A style, header-row style, data-row style).
    $styleSrc.Range("A2").Copy()
    $ws.Range("A1:A2").PasteSpecial(-4122)

    $hdrRange = $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item(1, $lastCol))
    $styleSrc.Range("B1").Copy()
    $hdrRange.PasteSpecial(-4122)

    $dataRange = $ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item(2, $lastCol))
    $styleSrc.Range("B2").Copy()
    $dataRange.PasteSpecial(-4122)
}

$styleTemplate = $wb.Worksheets.Item("土地")

# ---------------------------------------------------------------------------
# 1. Rename the existing "事業投資" sheet; its old content will be replaced later
#    (after we've copied it out into the new "事業投資" sheet created below).
# ---------------------------------------------------------------------------
$renamed = $wb.Worksheets.Item("事業投資")
$renamed.Name = "具有相當價值之財產"

# ---------------------------------------------------------------------------
# 2. Insert the three new sheets after it, in left-to-right order:
#    保險, 債權, 事業投資
# ---------------------------------------------------------------------------
$sInsurance = $wb.Worksheets.Add($null, $renamed)
$sInsurance.Name = "保險"

$sBond = $wb.Worksheets.Add($null, $sInsurance)
$sBond.Name = "債權"

$sBiz = $wb.Worksheets.Add($null, $sBond)
$sBiz.Name = "事業投資"

# ---------------------------------------------------------------------------
# 3. Populate the new "事業投資" sheet with the data that used to live in the
#    original "事業投資" sheet (now renamed to "具有相當價值之財產").
# ---------------------------------------------------------------------------
$sBiz.Cells.Clear()
$bizRows = @(
    @(199, "鄭汝芬", "任豐企業股份有限公司", "彰化縣西德里舜耕路47號", 1600000, "78年04月20日", "合資"),
    @(200, "鄭汝芬", "聯拳電信股份有限公司", "臺北市八德路2段232號4樓", 30460, "94年01月13日", "合資"),
    @(201, "鄭汝芬", "威寶電信股份有限公司", "臺北市瑞光路358巷36號5樓", 1222200, "94年02月02日", "合資"),
    @(202, "謝新隆", "任豐企業股份有限公司", "彰化縣西德里舜耕路*47號", 800000, "72年04月18曰", "合資"),
    @(203, "謝新隆", "埤頭液化煤氣行", "彰化縣復興路12號", 300000, "86年02月19日", "合夥"),
    @(204, "謝新隆", "豐宜實業股份有限公司", "高雄市崙北巷11號", 1000000, "91年10月23日", "合資"),
    @(205, "謝新隆", "六八煤氣行", "彰化縣慶平路83號", 120000, "74年04月20日", "合夥"),
    @(206, "謝新隆", "聯華電信股份有限公司", "臺北市八德路2段232號4樓", 134580, "94年01月13曰", "合"),
    @(207, "謝新隆", "亞太電信股份有限公司", "臺北市經貿二路66號12樓", 1000000, "97年12月31曰", "合資"),
    @(208, "謝新隆", "威寶電信股份肴限公司", "臺北市瑞光路358巷36號5樓", 726600, "94年02月02日", "合資"),
    @(209, "謝新隆", "全家福通訊科孩股份有限公司", "彰化縣中山路3段270巷11號1樓", 200000, "98年05月19日", "合資")
)

for ($i = 0; $i -lt $bizRows.Count; $i++) {
    $r = $i + 1
    $row = $bizRows[$i]
    $sBiz.Cells.Item($r, 1).Value = $row[0]
    $sBiz.Cells.Item($r, 2).Value = $row[1]
    $sBiz.Cells.Item($r, 3).Value = $row[2]
    $sBiz.Cells.Item($r, 4).Value = $row[3]
    $sBiz.Cells.Item($r, 5).Value = $row[4]
    $sBiz.Cells.Item($r, 6).Value = $row[5]
    $sBiz.Cells.Item($r, 7).Value = $row[6]
}

$styleTemplate.Range("A2").Copy()
$sBiz.Range("A1:A11").PasteSpecial(-4122)
$styleTemplate.Range("B2").Copy()
$sBiz.Range("B1:G11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Replace the (renamed) "具有相當價值之財產" sheet content with a single
#    jewelry ("珠寶") record.
# ---------------------------------------------------------------------------
$renamed.Cells.Clear()
$jewelHeaders = @("name","quantity","owner","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
$jewelValues  = @(150, "珠寶", 10, "鄭汝芬", 2150000, "otherbonds", "normal", "2011-11-21", "鄭汝芬", 1713, "tmpcd8e1", 150)
Write-HeaderAndRow $renamed $jewelHeaders $jewelValues $styleTemplate

# ---------------------------------------------------------------------------
# 5. Populate the "保險" (insurance) sheet.
# ---------------------------------------------------------------------------
$sInsurance.Cells.Clear()
$insHeaders = @("SW(罕發）勃迪","LBJ钿(帀召）**","盤m","7¥取茗Y揆褂","Y鄯勘","m#")
$insValues  = @(156, "SW(罕發）勃迪", "LBJ钿(帀召）**", "盤m", "7¥取茗Y揆褂", "Y鄯勘", "m#")
Write-HeaderAndRow $sInsurance $insHeaders $insValues $styleTemplate

# ---------------------------------------------------------------------------
# 6. Populate the "債權" (receivables) sheet.
# ---------------------------------------------------------------------------
$sBond.Cells.Clear()
$bondHeaders = @("OOOOSVZ","mm",1,"mm")
$bondValues  = @(184, "OOOOSVZ", "mm", "mm", 1)
Write-HeaderAndRow $sBond $bondHeaders $bondValues $styleTemplate

# ---------------------------------------------------------------------------
# 7. Rewrite the "其他有價證券" sheet (5th sheet) with a single, clean record.
# ---------------------------------------------------------------------------
$otherSec = $wb.Worksheets.Item("其他有價證券")
$otherSec.Cells.Clear()
$osHeaders = @("name","owner","quantity","face_value","currency","total","property_category","category","date","legislator_name","legislator_id","source_file","index")
$osValues  = @(146, "山隆通運股份有限公司", "鄭汝芬", 1, 10, "新臺幣", 10, "otherbonds", "normal", "2011-11-21", "鄭汝芬", 1713, "tmpcd8e1", 146)
Write-HeaderAndRow $otherSec $osHeaders $osValues $styleTemplate
